$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 272.0909
$ws.Range("I18").Value = 272.0909
$ws.Range("K18").Value = 272.0909
$ws.Range("M18").Value = 11.90910000000002

$ws.Range("H40").Value = 2180.5386
$ws.Range("I40").Value = 2136
$ws.Range("J40").Value = 2221.7778
$ws.Range("K40").Value = 2136
$ws.Range("L40").Value = 2221.7778
$ws.Range("M40").Value = -1961
$ws.Range("N40").Value = -2571.7778

$ws.Range("H43").Value = 421462.66
$ws.Range("I43").Value = 838391.2
$ws.Range("J43").Value = 4534.1665
$ws.Range("K43").Value = 838391.2
$ws.Range("L43").Value = 4534.1665
$ws.Range("M43").Value = -838322.2
$ws.Range("N43").Value = -4672.1665

$ws.Range("H45").Value = 665
$ws.Range("I45").Value = 665
$ws.Range("K45").Value = 1995
$ws.Range("M45").Value = -1803

$ws.Range("H48").Value = 5000
$ws.Range("J48").Value = 5000
$ws.Range("L48").Value = 15000
$ws.Range("N48").Value = -15584

$ws.Range("H55").Value = 299.4
$ws.Range("I55").Value = 299.25
$ws.Range("K55").Value = 299.25
$ws.Range("M55").Value = -85.25

$ws.Range("H56").Value = 5000
$ws.Range("J56").Value = 5000
$ws.Range("L56").Value = 15000
$ws.Range("N56").Value = -16068

$ws.Range("H74").Value = 8921.875
$ws.Range("I74").Value = 8315.799999999999
$ws.Range("K74").Value = 8315.799999999999
$ws.Range("M74").Value = -7379.799999999999

$ws.Range("H76").Value = 5002
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 5004
$ws.Range("K76").Value = 5000
$ws.Range("L76").Value = 5004
$ws.Range("M76").Value = -4685
$ws.Range("N76").Value = -5634

$ws.Range("H77").Value = 8921.875
$ws.Range("I77").Value = 8315.799999999999
$ws.Range("K77").Value = 41579
$ws.Range("M77").Value = -36899

$ws.Range("H79").Value = 5002
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 5004
$ws.Range("K79").Value = 5000
$ws.Range("L79").Value = 5004
$ws.Range("M79").Value = -3908
$ws.Range("N79").Value = -7188

$ws.Range("H82").Value = 1021.5
$ws.Range("I82").Value = 1021.5
$ws.Range("K82").Value = 3064.5
$ws.Range("M82").Value = -2658.5

$ws.Range("H85").Value = 1021.5
$ws.Range("I85").Value = 1021.5
$ws.Range("K85").Value = 3064.5
$ws.Range("M85").Value = -1660.5

$ws.Range("H86").Value = 2790.4348
$ws.Range("I86").Value = 2606.5625
$ws.Range("J86").Value = 3210.7144
$ws.Range("K86").Value = 2606.5625
$ws.Range("L86").Value = 3210.7144
$ws.Range("M86").Value = -1483.5625
$ws.Range("N86").Value = -5456.7144

$ws.Range("H89").Value = 2790.4348
$ws.Range("I89").Value = 2606.5625
$ws.Range("J89").Value = 3210.7144
$ws.Range("K89").Value = 13032.8125
$ws.Range("L89").Value = 16053.572
$ws.Range("M89").Value = -7416.8125
$ws.Range("N89").Value = -27285.572

$ws.Range("H98").Value = 1807.7428
$ws.Range("I98").Value = 1807.7428
$ws.Range("K98").Value = 1807.7428
$ws.Range("M98").Value = -309.7428

$ws.Range("H100").Value = 40169.58
$ws.Range("J100").Value = 5501.5
$ws.Range("L100").Value = 5501.5
$ws.Range("N100").Value = -6583.5

$ws.Range("H112").Value = 1869.6129
$ws.Range("J112").Value = 1929.6207
$ws.Range("L112").Value = 5788.8621
$ws.Range("N112").Value = -8004.8621

$ws.Range("H122").Value = 1807.7428
$ws.Range("I122").Value = 1807.7428
$ws.Range("K122").Value = 5423.2284
$ws.Range("M122").Value = -2973.2284

$ws.Range("H132").Value = 3247.8333
$ws.Range("I132").Value = 3171.6956
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 9515.086800000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -6985.086800000001
$ws.Range("N132").Value = -20057

$ws.Range("H135").Value = 3473
$ws.Range("I135").Value = 3692.75
$ws.Range("K135").Value = 33234.75
$ws.Range("M135").Value = -30699.75

$ws.Range("H137").Value = 10476.898
$ws.Range("I137").Value = 4670.0435
$ws.Range("K137").Value = 14010.1305
$ws.Range("M137").Value = -11460.1305

$ws.Range("H138").Value = 2786.3389
$ws.Range("I138").Value = 3459.9285
$ws.Range("J138").Value = 2576.7778
$ws.Range("K138").Value = 10379.7855
$ws.Range("L138").Value = 7730.3334
$ws.Range("M138").Value = -5239.7855
$ws.Range("N138").Value = -18010.3334

$ws.Range("H141").Value = 2776.5
$ws.Range("I141").Value = 2836.6875
$ws.Range("K141").Value = 8510.0625
$ws.Range("M141").Value = -3330.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6041.8906
$ws.Range("I32").Value = 6041.8906
$ws.Range("K32").Value = 6041.8906
$ws.Range("M32").Value = -5754.8906

$ws.Range("H41").Value = 3606.5
$ws.Range("I41").Value = 2982.1428
$ws.Range("J41").Value = 7977
$ws.Range("K41").Value = 2982.1428
$ws.Range("L41").Value = 7977
$ws.Range("M41").Value = -2568.1428
$ws.Range("N41").Value = -8805

$ws.Range("H45").Value = 13591
$ws.Range("I45").Value = 16205.143
$ws.Range("J45").Value = 4441.5
$ws.Range("K45").Value = 16205.143
$ws.Range("L45").Value = 4441.5
$ws.Range("M45").Value = -15828.143
$ws.Range("N45").Value = -5195.5

$ws.Range("H61").Value = 8411.6
$ws.Range("I61").Value = 5791.6665
$ws.Range("J61").Value = 14524.777
$ws.Range("K61").Value = 5791.6665
$ws.Range("L61").Value = 14524.777
$ws.Range("M61").Value = -5579.6665
$ws.Range("N61").Value = -14948.777

$ws.Range("H74").Value = 11044.195
$ws.Range("I74").Value = 11115.586
$ws.Range("K74").Value = 11115.586
$ws.Range("M74").Value = -10241.586

$ws.Range("H77").Value = 11044.195
$ws.Range("I77").Value = 11115.586
$ws.Range("K77").Value = 55577.92999999999
$ws.Range("M77").Value = -51209.92999999999

$ws.Range("H88").Value = 1673.591
$ws.Range("I88").Value = 1247.1666
$ws.Range("J88").Value = 1833.5
$ws.Range("K88").Value = 1247.1666
$ws.Range("L88").Value = 1833.5
$ws.Range("M88").Value = -841.1666
$ws.Range("N88").Value = -2645.5

$ws.Range("H91").Value = 1673.591
$ws.Range("I91").Value = 1247.1666
$ws.Range("J91").Value = 1833.5
$ws.Range("K91").Value = 1247.1666
$ws.Range("L91").Value = 1833.5
$ws.Range("M91").Value = 156.8334
$ws.Range("N91").Value = -4641.5

$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H122").Value = 2511.75
$ws.Range("I122").Value = 2616.647
$ws.Range("J122").Value = 2257
$ws.Range("K122").Value = 7849.941
$ws.Range("L122").Value = 6771
$ws.Range("M122").Value = -5399.941
$ws.Range("N122").Value = -11671

$ws.Range("H132").Value = 2932.6572
$ws.Range("I132").Value = 3119.3928
$ws.Range("K132").Value = 9358.178400000001
$ws.Range("M132").Value = -6828.178400000001

$ws.Range("H136").Value = 8411.6
$ws.Range("I136").Value = 5791.6665
$ws.Range("J136").Value = 14524.777
$ws.Range("K136").Value = 17374.9995
$ws.Range("L136").Value = 43574.331
$ws.Range("M136").Value = -14824.9995
$ws.Range("N136").Value = -48674.331

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 8014
$ws.Range("I25").Value = 8014
$ws.Range("K25").Value = 8014
$ws.Range("M25").Value = -7779

$ws.Range("H82").Value = 56045.277
$ws.Range("J82").Value = 90994.60000000001
$ws.Range("L82").Value = 90994.60000000001
$ws.Range("N82").Value = -91760.60000000001

$ws.Range("H85").Value = 56045.277
$ws.Range("J85").Value = 90994.60000000001
$ws.Range("L85").Value = 90994.60000000001
$ws.Range("N85").Value = -93646.60000000001

$ws.Range("H86").Value = 530769.2
$ws.Range("I86").Value = 2003185.6
$ws.Range("J86").Value = 4906.2144
$ws.Range("K86").Value = 2003185.6
$ws.Range("L86").Value = 4906.2144
$ws.Range("M86").Value = -2002062.6
$ws.Range("N86").Value = -7152.2144

$ws.Range("H89").Value = 530769.2
$ws.Range("I89").Value = 2003185.6
$ws.Range("J89").Value = 4906.2144
$ws.Range("K89").Value = 10015928
$ws.Range("L89").Value = 24531.072
$ws.Range("M89").Value = -10010312
$ws.Range("N89").Value = -35763.072

$ws.Range("H105").Value = 4847.0386
$ws.Range("I105").Value = 4595
$ws.Range("J105").Value = 5323.1113
$ws.Range("K105").Value = 4595
$ws.Range("L105").Value = 5323.1113
$ws.Range("M105").Value = -2848
$ws.Range("N105").Value = -8817.1113

$ws.Range("H107").Value = 9714.143
$ws.Range("I107").Value = 13747.25
$ws.Range("J107").Value = 4336.6665
$ws.Range("K107").Value = 13747.25
$ws.Range("L107").Value = 4336.6665
$ws.Range("M107").Value = -11827.25
$ws.Range("N107").Value = -8176.6665

$ws.Range("H134").Value = 8130.795
$ws.Range("I134").Value = 3503.5938
$ws.Range("K134").Value = 10510.7814
$ws.Range("M134").Value = -7975.7814

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 19608.4
$ws.Range("J28").Value = 19608.4
$ws.Range("L28").Value = 19608.4
$ws.Range("N28").Value = -20098.4

$ws.Range("H31").Value = 3658
$ws.Range("I31").Value = 2788.1
$ws.Range("K31").Value = 2788.1
$ws.Range("M31").Value = -2493.1

$ws.Range("H34").Value = 3658
$ws.Range("I34").Value = 2788.1
$ws.Range("K34").Value = 2788.1
$ws.Range("M34").Value = -2586.1

$ws.Range("H43").Value = 15951.833
$ws.Range("J43").Value = 15951.833
$ws.Range("L43").Value = 15951.833
$ws.Range("N43").Value = -16319.833

$ws.Range("H101").Value = 15951.833
$ws.Range("J101").Value = 15951.833
$ws.Range("L101").Value = 15951.833
$ws.Range("N101").Value = -22441.833

$ws.Range("H132").Value = 27657.611
$ws.Range("I132").Value = 21613.576
$ws.Range("J132").Value = 33269.93
$ws.Range("K132").Value = 64840.728
$ws.Range("L132").Value = 99809.79000000001
$ws.Range("M132").Value = -62310.728
$ws.Range("N132").Value = -104869.79

$ws.Range("H133").Value = 600356
$ws.Range("I133").Value = 65998.39999999999
$ws.Range("K133").Value = 65998.39999999999
$ws.Range("M133").Value = -63468.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 48.473682
$ws.Range("J2").Value = 54.666668
$ws.Range("L2").Value = 328.000008
$ws.Range("N2").Value = -554.000008

$ws.Range("H5").Value = 2267.5898
$ws.Range("I5").Value = 1171.25
$ws.Range("K5").Value = 3513.75
$ws.Range("M5").Value = -3401.75

$ws.Range("H11").Value = 123878.97
$ws.Range("I11").Value = 289.97702
$ws.Range("K11").Value = 869.9310599999999
$ws.Range("M11").Value = -729.9310599999999

$ws.Range("H34").Value = 3273.8096
$ws.Range("I34").Value = 299.44446
$ws.Range("J34").Value = 5504.5835
$ws.Range("K34").Value = 898.33338
$ws.Range("L34").Value = 16513.7505
$ws.Range("M34").Value = -814.33338
$ws.Range("N34").Value = -16681.7505

$ws.Range("H39").Value = 4022.1667
$ws.Range("J39").Value = 5324.9165
$ws.Range("L39").Value = 15974.7495
$ws.Range("N39").Value = -16562.7495

$ws.Range("H55").Value = 1199.8889
$ws.Range("I55").Value = 121.28571
$ws.Range("J55").Value = 4975
$ws.Range("K55").Value = 363.85713
$ws.Range("L55").Value = 14925
$ws.Range("M55").Value = -186.85713
$ws.Range("N55").Value = -15279

$ws.Range("H63").Value = 6799.6
$ws.Range("I63").Value = 7999.3335
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 23998.0005
$ws.Range("L63").Value = 15000
$ws.Range("M63").Value = -23249.0005
$ws.Range("N63").Value = -16498

$ws.Range("H64").Value = 1360.75
$ws.Range("I64").Value = 1249.5
$ws.Range("J64").Value = 1472
$ws.Range("K64").Value = 3748.5
$ws.Range("L64").Value = 4416
$ws.Range("M64").Value = -3478.5
$ws.Range("N64").Value = -4956

$ws.Range("H66").Value = 6799.6
$ws.Range("I66").Value = 7999.3335
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 71994.0015
$ws.Range("L66").Value = 45000
$ws.Range("M66").Value = -68250.0015
$ws.Range("N66").Value = -52488

$ws.Range("H67").Value = 1360.75
$ws.Range("I67").Value = 1249.5
$ws.Range("J67").Value = 1472
$ws.Range("K67").Value = 3748.5
$ws.Range("L67").Value = 4416
$ws.Range("M67").Value = -2812.5
$ws.Range("N67").Value = -6288

$ws.Range("H70").Value = 176248.75
$ws.Range("I70").Value = 176248.75
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 528746.25
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -528431.25
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 176248.75
$ws.Range("I73").Value = 176248.75
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 528746.25
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -527654.25
$ws.Range("N73").ClearContents()

$ws.Range("H74").Value = 3999999
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 3999999
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H92").Value = 1695
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 1695
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 5085
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -7581

$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()

$ws.Range("H103").Value = 1441.4445
$ws.Range("I103").Value = 1492.2858
$ws.Range("J103").Value = 1263.5
$ws.Range("K103").Value = 4476.857400000001
$ws.Range("L103").Value = 3790.5
$ws.Range("M103").Value = -3597.857400000001
$ws.Range("N103").Value = -5548.5

$ws.Range("H104").Value = 10981.793
$ws.Range("I104").Value = 7042
$ws.Range("J104").Value = 12119.955
$ws.Range("K104").Value = 21126
$ws.Range("L104").Value = 36359.865
$ws.Range("M104").Value = -18505
$ws.Range("N104").Value = -41601.865

$ws.Range("H109").Value = 1726.3334
$ws.Range("I109").Value = 1726.3334
$ws.Range("K109").Value = 5179.0002
$ws.Range("M109").Value = -4139.0002

$ws.Range("H117").Value = 1771.2222
$ws.Range("I117").Value = 791.5714
$ws.Range("J117").Value = 5200
$ws.Range("K117").Value = 2374.7142
$ws.Range("L117").Value = 15600
$ws.Range("M117").Value = 1067.2858
$ws.Range("N117").Value = -22484

$ws.Range("H131").Value = 2907.8857
$ws.Range("I131").Value = 1079.1875
$ws.Range("J131").Value = 4447.8423
$ws.Range("K131").Value = 3237.5625
$ws.Range("L131").Value = 13343.5269
$ws.Range("M131").Value = 1802.4375
$ws.Range("N131").Value = -23423.5269

$ws.Range("H132").Value = 4764726
$ws.Range("I132").Value = 3026.8462
$ws.Range("J132").Value = 12502488
$ws.Range("K132").Value = 27241.6158
$ws.Range("L132").Value = 112522392
$ws.Range("M132").Value = -24711.6158
$ws.Range("N132").Value = -112527452

$ws.Range("H135").Value = 2267.5898
$ws.Range("I135").Value = 1171.25
$ws.Range("K135").Value = 10541.25
$ws.Range("M135").Value = -8006.25

$ws.Range("H137").Value = 7968.5
$ws.Range("I137").Value = 8048.25
$ws.Range("K137").Value = 24144.75
$ws.Range("M137").Value = -19044.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3179.742
$ws.Range("I122").Value = 3392.25
$ws.Range("K122").Value = 10176.75
$ws.Range("M122").Value = -7726.75

$ws.Range("H132").Value = 7721.826
$ws.Range("I132").Value = 6692.316
$ws.Range("J132").Value = 12612
$ws.Range("K132").Value = 20076.948
$ws.Range("L132").Value = 37836
$ws.Range("M132").Value = -17546.948
$ws.Range("N132").Value = -42896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

$ws.Range("H22").Value = 899.8333
$ws.Range("I22").Value = 689.5
$ws.Range("K22").Value = 689.5
$ws.Range("M22").Value = -394.5

$ws.Range("H27").Value = 899.8333
$ws.Range("I27").Value = 689.5
$ws.Range("K27").Value = 689.5
$ws.Range("M27").Value = -582.5

$ws.Range("H40").Value = 4556.2607
$ws.Range("I40").Value = 4500.067
$ws.Range("J40").Value = 4661.625
$ws.Range("K40").Value = 4500.067
$ws.Range("L40").Value = 4661.625
$ws.Range("M40").Value = -4364.067
$ws.Range("N40").Value = -4933.625

$ws.Range("H46").Value = 1740.1
$ws.Range("I46").Value = 1080.9
$ws.Range("J46").Value = 2069.7
$ws.Range("K46").Value = 1080.9
$ws.Range("L46").Value = 2069.7
$ws.Range("M46").Value = -892.9000000000001
$ws.Range("N46").Value = -2445.7

$ws.Range("H100").Value = 3000
$ws.Range("I100").Value = 3000
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459

$ws.Range("H102").Value = 99999
$ws.Range("J102").Value = 99999
$ws.Range("L102").Value = 99999
$ws.Range("N102").Value = -106489

$ws.Range("H103").Value = 36301
$ws.Range("J103").Value = 36301
$ws.Range("L103").Value = 36301
$ws.Range("N103").Value = -38645

$ws.Range("H110").Value = 39999
$ws.Range("J110").Value = 39999
$ws.Range("L110").Value = 39999
$ws.Range("N110").Value = -48179

$ws.Range("H136").Value = 5262.5293
$ws.Range("I136").Value = 4681.8364
$ws.Range("K136").Value = 14045.5092
$ws.Range("M136").Value = -11495.5092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2342.6843
$ws.Range("J81").Value = 3000
$ws.Range("L81").Value = 6000
$ws.Range("N81").Value = -8122

$ws.Range("H84").Value = 2342.6843
$ws.Range("J84").Value = 3000
$ws.Range("L84").Value = 30000
$ws.Range("N84").Value = -40608

$ws.Range("H103").Value = 26880.334
$ws.Range("J103").Value = 26880.334
$ws.Range("L103").Value = 26880.334
$ws.Range("N103").Value = -29224.334

$ws.Range("H107").Value = 1662.3889
$ws.Range("J107").Value = 2956.2
$ws.Range("L107").Value = 8868.599999999999
$ws.Range("N107").Value = -12708.6

$ws.Range("H119").Value = 415000
$ws.Range("J119").Value = 415000
$ws.Range("L119").Value = 415000
$ws.Range("N119").Value = -424676

$ws.Range("H122").Value = 3659.7646
$ws.Range("I122").Value = 2949
$ws.Range("J122").Value = 4459.375
$ws.Range("K122").Value = 8847
$ws.Range("L122").Value = 13378.125
$ws.Range("M122").Value = -6397
$ws.Range("N122").Value = -18278.125

$ws.Range("H123").Value = 49999.91
$ws.Range("J123").Value = 49999.91
$ws.Range("L123").Value = 49999.91
$ws.Range("N123").Value = -59799.91

$ws.Range("H132").Value = 146802
$ws.Range("I132").Value = 202342.02
$ws.Range("J132").Value = 29057.16
$ws.Range("K132").Value = 607026.0599999999
$ws.Range("L132").Value = 87171.48
$ws.Range("M132").Value = -604496.0599999999
$ws.Range("N132").Value = -92231.48
